# Refresh the charging-station "idle since last charge" report: bump the
# "as of" timestamp (column D) and recompute idle hours (column E) for every
# existing row, swap in several terminals' updated charge-end times/terminal
# ids (columns A-C) to reflect newly completed sessions, and append four new
# rows (35-38) for terminals that have since gone idle.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = 45934.358252314814
$ws.Cells.Item(2, 5).Value = 165.62194444437046

# Row 3
$ws.Cells.Item(3, 4).Value = 45934.358252314814
$ws.Cells.Item(3, 5).Value = 98.885000000009313

# Row 4
$ws.Cells.Item(4, 4).Value = 45934.358252314814
$ws.Cells.Item(4, 5).Value = 92.188611111079808

# Row 5
$ws.Cells.Item(5, 4).Value = 45934.358252314814
$ws.Cells.Item(5, 5).Value = 89.753611111140344

# Row 6
$ws.Cells.Item(6, 4).Value = 45934.358252314814
$ws.Cells.Item(6, 5).Value = 75.16777777770767

# Row 7
$ws.Cells.Item(7, 4).Value = 45934.358252314814
$ws.Cells.Item(7, 5).Value = 66.784999999974389

# Row 8
$ws.Cells.Item(8, 1).Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Cells.Item(8, 2).Value = "102号直流"
$ws.Cells.Item(8, 3).Value = 45931.647812499999
$ws.Cells.Item(8, 4).Value = 45934.358252314814
$ws.Cells.Item(8, 5).Value = 65.050555555557366

# Row 9
$ws.Cells.Item(9, 2).Value = "801号直流"
$ws.Cells.Item(9, 3).Value = 45932.072800925926
$ws.Cells.Item(9, 4).Value = 45934.358252314814
$ws.Cells.Item(9, 5).Value = 54.850833333330229

# Row 10
$ws.Cells.Item(10, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(10, 2).Value = "603号直流"
$ws.Cells.Item(10, 3).Value = 45932.081099537034
$ws.Cells.Item(10, 4).Value = 45934.358252314814
$ws.Cells.Item(10, 5).Value = 54.651666666730307

# Row 11
$ws.Cells.Item(11, 2).Value = "905号直流"
$ws.Cells.Item(11, 3).Value = 45932.25236111111
$ws.Cells.Item(11, 4).Value = 45934.358252314814
$ws.Cells.Item(11, 5).Value = 50.541388888901565

# Row 12
$ws.Cells.Item(12, 1).Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Cells.Item(12, 2).Value = "402号直流"
$ws.Cells.Item(12, 3).Value = 45932.255543981482
$ws.Cells.Item(12, 4).Value = 45934.358252314814
$ws.Cells.Item(12, 5).Value = 50.464999999967404

# Row 13
$ws.Cells.Item(13, 1).Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Cells.Item(13, 2).Value = "209号直流"
$ws.Cells.Item(13, 3).Value = 45933.030578703707
$ws.Cells.Item(13, 4).Value = 45934.358252314814
$ws.Cells.Item(13, 5).Value = 31.864166666578967

# Row 14
$ws.Cells.Item(14, 2).Value = "805号直流"
$ws.Cells.Item(14, 3).Value = 45933.035462962966
$ws.Cells.Item(14, 4).Value = 45934.358252314814
$ws.Cells.Item(14, 5).Value = 31.746944444370456

# Row 15
$ws.Cells.Item(15, 2).Value = "406号直流"
$ws.Cells.Item(15, 3).Value = 45933.039143518516
$ws.Cells.Item(15, 4).Value = 45934.358252314814
$ws.Cells.Item(15, 5).Value = 31.658611111168284

# Row 16
$ws.Cells.Item(16, 1).Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Cells.Item(16, 2).Value = "109号直流"
$ws.Cells.Item(16, 3).Value = 45933.063738425924
$ws.Cells.Item(16, 4).Value = 45934.358252314814
$ws.Cells.Item(16, 5).Value = 31.068333333358169

# Row 17
$ws.Cells.Item(17, 2).Value = "002A号直流"
$ws.Cells.Item(17, 3).Value = 45933.065972222219
$ws.Cells.Item(17, 4).Value = 45934.358252314814
$ws.Cells.Item(17, 5).Value = 31.014722222287674

# Row 18
$ws.Cells.Item(18, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(18, 2).Value = "904号直流"
$ws.Cells.Item(18, 3).Value = 45933.088784722226
$ws.Cells.Item(18, 4).Value = 45934.358252314814
$ws.Cells.Item(18, 5).Value = 30.46722222212702

# Row 19
$ws.Cells.Item(19, 1).Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Cells.Item(19, 2).Value = "104号直流"
$ws.Cells.Item(19, 3).Value = 45933.200601851851
$ws.Cells.Item(19, 4).Value = 45934.358252314814
$ws.Cells.Item(19, 5).Value = 27.783611111110076

# Row 20
$ws.Cells.Item(20, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(20, 2).Value = "A01号直流"
$ws.Cells.Item(20, 3).Value = 45933.222430555557
$ws.Cells.Item(20, 4).Value = 45934.358252314814
$ws.Cells.Item(20, 5).Value = 27.259722222166602

# Row 21
$ws.Cells.Item(21, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(21, 2).Value = "504号直流"
$ws.Cells.Item(21, 3).Value = 45933.270925925928
$ws.Cells.Item(21, 4).Value = 45934.358252314814
$ws.Cells.Item(21, 5).Value = 26.095833333267365

# Row 22
$ws.Cells.Item(22, 2).Value = "103号直流"
$ws.Cells.Item(22, 3).Value = 45933.305023148147
$ws.Cells.Item(22, 4).Value = 45934.358252314814
$ws.Cells.Item(22, 5).Value = 25.277500000025611

# Row 23
$ws.Cells.Item(23, 1).Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Cells.Item(23, 2).Value = "107号直流"
$ws.Cells.Item(23, 3).Value = 45933.386261574073
$ws.Cells.Item(23, 4).Value = 45934.358252314814
$ws.Cells.Item(23, 5).Value = 23.327777777798474

# Row 24
$ws.Cells.Item(24, 1).Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Cells.Item(24, 2).Value = "205号直流"
$ws.Cells.Item(24, 3).Value = 45933.405219907407
$ws.Cells.Item(24, 4).Value = 45934.358252314814
$ws.Cells.Item(24, 5).Value = 22.872777777782176

# Row 25
$ws.Cells.Item(25, 1).Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Cells.Item(25, 2).Value = "207号直流"
$ws.Cells.Item(25, 3).Value = 45933.481249999997
$ws.Cells.Item(25, 4).Value = 45934.358252314814
$ws.Cells.Item(25, 5).Value = 21.048055555613246

# Row 26
$ws.Cells.Item(26, 1).Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Cells.Item(26, 2).Value = "002B号直流"
$ws.Cells.Item(26, 3).Value = 45933.517557870371
$ws.Cells.Item(26, 4).Value = 45934.358252314814
$ws.Cells.Item(26, 5).Value = 20.176666666637175

# Row 27
$ws.Cells.Item(27, 2).Value = "105号直流"
$ws.Cells.Item(27, 3).Value = 45933.521967592591
$ws.Cells.Item(27, 4).Value = 45934.358252314814
$ws.Cells.Item(27, 5).Value = 20.070833333360497

# Row 28
$ws.Cells.Item(28, 2).Value = "306号直流"
$ws.Cells.Item(28, 3).Value = 45933.529745370368
$ws.Cells.Item(28, 4).Value = 45934.358252314814
$ws.Cells.Item(28, 5).Value = 19.884166666714009

# Row 29
$ws.Cells.Item(29, 2).Value = "705号直流"
$ws.Cells.Item(29, 3).Value = 45933.537893518522
$ws.Cells.Item(29, 4).Value = 45934.358252314814
$ws.Cells.Item(29, 5).Value = 19.688611111021601

# Row 30
$ws.Cells.Item(30, 1).Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Cells.Item(30, 2).Value = "011A号直流"
$ws.Cells.Item(30, 3).Value = 45933.552754629629
$ws.Cells.Item(30, 4).Value = 45934.358252314814
$ws.Cells.Item(30, 5).Value = 19.331944444449618

# Row 31
$ws.Cells.Item(31, 2).Value = "703号直流"
$ws.Cells.Item(31, 3).Value = 45933.557268518518
$ws.Cells.Item(31, 4).Value = 45934.358252314814
$ws.Cells.Item(31, 5).Value = 19.223611111112405

# Row 32
$ws.Cells.Item(32, 1).Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Cells.Item(32, 2).Value = "203号直流"
$ws.Cells.Item(32, 3).Value = 45933.64603009259
$ws.Cells.Item(32, 4).Value = 45934.358252314814
$ws.Cells.Item(32, 5).Value = 17.093333333381452

# Row 33
$ws.Cells.Item(33, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(33, 2).Value = "601号直流"
$ws.Cells.Item(33, 3).Value = 45933.663263888891
$ws.Cells.Item(33, 4).Value = 45934.358252314814
$ws.Cells.Item(33, 5).Value = 16.679722222150303

# Row 34
$ws.Cells.Item(34, 2).Value = "104号直流"
$ws.Cells.Item(34, 3).Value = 45933.714675925927
$ws.Cells.Item(34, 4).Value = 45934.358252314814
$ws.Cells.Item(34, 5).Value = 15.445833333302289

# Row 35
$ws.Cells.Item(35, 1).Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Cells.Item(35, 2).Value = "404号直流"
$ws.Cells.Item(35, 3).Value = 45933.724224537036
$ws.Cells.Item(35, 4).Value = 45934.358252314814
$ws.Cells.Item(35, 5).Value = 15.216666666674428

# Row 36
$ws.Cells.Item(36, 1).Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Cells.Item(36, 2).Value = "906号直流"
$ws.Cells.Item(36, 3).Value = 45933.791134259256
$ws.Cells.Item(36, 4).Value = 45934.358252314814
$ws.Cells.Item(36, 5).Value = 13.61083333339775

# Row 37
$ws.Cells.Item(37, 1).Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Cells.Item(37, 2).Value = "212号直流"
$ws.Cells.Item(37, 3).Value = 45933.821712962963
$ws.Cells.Item(37, 4).Value = 45934.358252314814
$ws.Cells.Item(37, 5).Value = 12.87694444443332

# Row 38
$ws.Cells.Item(38, 1).Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Cells.Item(38, 2).Value = "401号直流"
$ws.Cells.Item(38, 3).Value = 45933.847962962966
$ws.Cells.Item(38, 4).Value = 45934.358252314814
$ws.Cells.Item(38, 5).Value = 12.246944444370456

# Move the active selection to match the author's last cursor position.
$ws.Range("I9").Select()
